$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores these figures as literal text (inline strings),
# e.g. "319.34" or "3.30%", not as numbers/percentages. Force text storage
# by setting NumberFormat "@" before assignment, then restore the default
# "Normal" style afterwards so we do not leave a stray text-format style on
# the cell (matches the unformatted cells in the source).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "319.34"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "3.30%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "41.40"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "1.21%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.265"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "2.79%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07743"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1.60%"
$c.Style = "Normal"

$c = $ws.Range("B6")
$c.NumberFormat = "@"
$c.Value = "FTXToken"
$c.Style = "Normal"

$c = $ws.Range("C6")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.757"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "9.61%"
$c.Style = "Normal"

$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c.Style = "Normal"

$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.9439"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "4.00%"
$c.Style = "Normal"

$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "BTSEToken"
$c.Style = "Normal"

$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "2.425"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "-1.58%"
$c.Style = "Normal"

$c = $ws.Range("B9")
$c.NumberFormat = "@"
$c.Value = "LiechtensteinCryptoassetsExchange"
$c.Style = "Normal"

$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1246"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "-1.74%"
$c.Style = "Normal"

$c = $ws.Range("B10")
$c.NumberFormat = "@"
$c.Value = "WazirX"
$c.Style = "Normal"

$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1887"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "4.56%"
$c.Style = "Normal"

$c = $ws.Range("B11")
$c.NumberFormat = "@"
$c.Value = "MandalaExchangeToken"
$c.Style = "Normal"

$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.09242"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "2.19%"
$c.Style = "Normal"

$c = $ws.Range("B12")
$c.NumberFormat = "@"
$c.Value = "BitrueCoin"
$c.Style = "Normal"

$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.04306"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "0.18%"
$c.Style = "Normal"

$c = $ws.Range("B13")
$c.NumberFormat = "@"
$c.Value = "BitMartToken"
$c.Style = "Normal"

$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.1051"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "0.64%"
$c.Style = "Normal"

$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = "BitForexToken"
$c.Style = "Normal"

$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001281"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1.71%"
$c.Style = "Normal"

$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "TigerCash"
$c.Style = "Normal"

$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.005851"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "2.40%"
$c.Style = "Normal"

$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "UpBots"
$c.Style = "Normal"

$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.007491"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1,897.31%"
$c.Style = "Normal"

$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "LEO"
$c.Style = "Normal"

$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.342"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-0.25%"
$c.Style = "Normal"

$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = "GateToken"
$c.Style = "Normal"

$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.339"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "1.34%"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.31%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.793"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "12.80%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.1343"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-3.61%"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "5.10%"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-0.24%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.001268"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-0.24%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004126"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "2.03%"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "0.00%"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02557"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "5.71%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05332"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "2.14%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.007778"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "-0.73%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1316"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.36%"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.007043"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "3.47%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.001991"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "2.98%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.008265"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "12.62%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.3175"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-5.56%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006677"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-3.05%"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000750"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-0.13%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.2003"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "56.09%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.004201"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "39.92%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002101"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.13%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002001"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "-0.13%"
$c.Style = "Normal"

